# Generate Report for Archive
#
# 1. Update status text "Ready for handoff" -> "In Translation" everywhere
#    it appears (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "status" columns that held that text:
#      Overview: columns E (5) and F (6)
#      zh-cn:    column C (3)
#      de-de:    column C (3)
#    from ~17.22 chars wide down to ~13.41 chars wide.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# -- Update the status text --------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# -- Resize the status columns ------------------------------------------
# Target stored column width is 13.4101845877511 characters; the
# ColumnWidth property below is the input value that yields the closest
# achievable stored width.
$newColumnWidth = 12.576851254417766

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
